$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.906.63'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.216.56'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '241.01'
$ws.Range("D5").ClearFormats()
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '73.44'
$ws.Range("D7").ClearFormats()
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.607'
$ws.Range("D9").ClearFormats()
$ws.Range('E9').Value = '  -2.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '42.64'
$ws.Range("D10").ClearFormats()
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '7.09'
$ws.Range("D12").ClearFormats()
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').Value = '2.550.05'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '14.28'
$ws.Range("D15").ClearFormats()
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '0.837'
$ws.Range("D16").ClearFormats()
$ws.Range('E16').Value = '  -2.17%  '
$ws.Range('D17').Value = '2.220.31'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '41.892.09'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('E19').Value = '  +8.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '73.07'
$ws.Range("D20").ClearFormats()
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '6.16'
$ws.Range("D21").ClearFormats()
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '10.80'
$ws.Range("D22").ClearFormats()
$ws.Range('E22').Value = '  +21.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '229.42'
$ws.Range("D23").ClearFormats()
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('E24').Value = '  -6.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '11.70'
$ws.Range("D25").ClearFormats()
$ws.Range('E25').Value = '  +2.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  -1.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '2.20'
$ws.Range("D29").ClearFormats()
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '166.99'
$ws.Range("D30").ClearFormats()
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '5.65'
$ws.Range("D32").ClearFormats()
$ws.Range('E32').Value = '  +7.33%  '
$ws.Range('E33').Value = '  -3.96%  '
$ws.Range('E34').Value = '  -0.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '29.28'
$ws.Range("D35").ClearFormats()
$ws.Range('E35').Value = '  -6.22%  '
$ws.Range('E36').Value = '  -11.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '4.27'
$ws.Range("D37").ClearFormats()
$ws.Range('E37').Value = '  -5.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '0.0299'
$ws.Range("D38").ClearFormats()
$ws.Range('E38').Value = '  -5.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '13.76'
$ws.Range("D39").ClearFormats()
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '65.70'
$ws.Range("D40").ClearFormats()
$ws.Range('E40').Value = '  +8.14%  '
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '5.62'
$ws.Range("D42").ClearFormats()
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '8.70'
$ws.Range("D44").ClearFormats()
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '104.18'
$ws.Range("D45").ClearFormats()
$ws.Range('E45').Value = '  -2.96%  '
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '1.16'
$ws.Range("D49").ClearFormats()
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '2.69'
$ws.Range("D50").ClearFormats()
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('D51').Value = '2.424.30'
$ws.Range('E51').Value = '  -1.44%  '
